$d = $word.ActiveDocument

# Remove the leading "section number" + tab that precedes the heading text
# (e.g. "1\tSupplemental Results" -> "Supplemental Results") for every
# Heading 1-4 styled paragraph that starts with the auto-generated
# SectionNumber run followed by a tab character.

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Range.ParagraphStyle.NameLocal
    if ($styleName -eq "Heading 1" -or $styleName -eq "Heading 2" -or $styleName -eq "Heading 3" -or $styleName -eq "Heading 4") {
        $full = $p.Range.Text
        $tabIndex = $full.IndexOf([char]9)
        # The generated section number (e.g. "1.1.1.1") is short, so only treat
        # a tab found near the very start of the paragraph as the marker to strip.
        if ($tabIndex -ge 0 -and $tabIndex -le 10) {
            $start = $p.Range.Start
            $prefix = $d.Range($start, $start + $tabIndex + 1)
            $prefix.Text = ""
        }
    }
}
